$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Columns")

# New rows of unit-related analytes with Tab = "FieldResult"
$newRows = @(
    @("WaterTemperatureUnit", "FieldResult"),
    @("AirTemperatureUnit", "FieldResult"),
    @("AirWindSpeedUnit", "FieldResult"),
    @("WaterDOSatUnit", "FieldResult"),
    @("WaterSpConductivityUnit", "FieldResult"),
    @("WaterSalinityUnit", "FieldResult")
)

$startRow = 68
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 4).Value = $newRows[$i][1]
}

# Update the view: scroll to B59 and select D73
$ws.Activate()
$lastRow = $startRow + $newRows.Count - 1
$ws.Application.Goto($ws.Range("B59"), $false)
$ws.Range("D$lastRow").Select()
